$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cryptos list refresh (prices / 1h volume %) as of the GitHub Actions run.
# Rows 9/10 and 40/41 also swap coin identity (ranking reorder), so B/C/D/E
# are all rewritten for those two row pairs.

$ws.Range("D2").Value = "'28.260.59"
$ws.Range("E2").Value = "'  +5.08%  "

$ws.Range("D3").Value = "'1.799.37"
$ws.Range("E3").Value = "'  +3.51%  "

$ws.Range("D4").Value = "'1.000"
$ws.Range("E4").Value = "'  -0.01%  "

$ws.Range("D5").Value = "'316.45"
$ws.Range("E5").Value = "'  +1.87%  "

$ws.Range("D6").Value = "'1.000"
$ws.Range("E6").Value = "'  -0.03%  "

$ws.Range("D7").Value = "'0.5518"
$ws.Range("E7").Value = "'  +10.04%  "

$ws.Range("D8").Value = "'0.3877"
$ws.Range("E8").Value = "'  +8.57%  "

$ws.Range("B9").Value = "'Dogecoin"
$ws.Range("C9").Value = "'https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge"
$ws.Range("D9").Value = "'0.07577"
$ws.Range("E9").Value = "'  +4.56%  "

$ws.Range("B10").Value = "'OKB"
$ws.Range("C10").Value = "'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D10").Value = "'42.79"
$ws.Range("E10").Value = "'  +1.50%  "

$ws.Range("D11").Value = "'1.117"
$ws.Range("E11").Value = "'  +5.41%  "

$ws.Range("E12").Value = "'  -0.06%  "

$ws.Range("D13").Value = "'21.15"
$ws.Range("E13").Value = "'  +4.94%  "

$ws.Range("D14").Value = "'6.204"
$ws.Range("E14").Value = "'  +4.54%  "

$ws.Range("D15").Value = "'7.343"
$ws.Range("E15").Value = "'  +7.83%  "

$ws.Range("D16").Value = "'1.799.37"
$ws.Range("E16").Value = "'  +3.85%  "

$ws.Range("D17").Value = "'91.94"
$ws.Range("E17").Value = "'  +6.34%  "

$ws.Range("E18").Value = "'  +3.48%  "

$ws.Range("D19").Value = "'0.06453"
$ws.Range("E19").Value = "'  +0.64%  "

$ws.Range("D20").Value = "'1.000"
$ws.Range("E20").Value = "'  -0.08%  "

$ws.Range("E21").Value = "'  +4.60%  "

$ws.Range("D22").Value = "'5.969"
$ws.Range("E22").Value = "'  +4.11%  "

$ws.Range("D23").Value = "'28.275.18"
$ws.Range("E23").Value = "'  +4.88%  "

$ws.Range("D24").Value = "'11.50"
$ws.Range("E24").Value = "'  +2.11%  "

$ws.Range("D25").Value = "'2.154"
$ws.Range("E25").Value = "'  +5.28%  "

$ws.Range("D26").Value = "'157.92"
$ws.Range("E26").Value = "'  +2.82%  "

$ws.Range("D27").Value = "'20.67"
$ws.Range("E27").Value = "'  +4.62%  "

$ws.Range("D28").Value = "'2.398"
$ws.Range("E28").Value = "'  +8.28%  "

$ws.Range("D29").Value = "'2.004.31"
$ws.Range("E29").Value = "'  +1.52%  "

$ws.Range("D30").Value = "'123.54"
$ws.Range("E30").Value = "'  +3.35%  "

$ws.Range("D31").Value = "'1.132"
$ws.Range("E31").Value = "'  +8.79%  "

$ws.Range("D32").Value = "'0.1022"
$ws.Range("E32").Value = "'  +7.22%  "

$ws.Range("D33").Value = "'5.736"
$ws.Range("E33").Value = "'  +7.18%  "

$ws.Range("D34").Value = "'3.670"
$ws.Range("E34").Value = "'  +2.50%  "

$ws.Range("D35").Value = "'0.2330"
$ws.Range("E35").Value = "'  +16.44%  "

$ws.Range("D36").Value = "'0.06378"
$ws.Range("E36").Value = "'  +8.42%  "

$ws.Range("D37").Value = "'8.851"
$ws.Range("E37").Value = "'  +16.11%  "

$ws.Range("D38").Value = "'0.02313"
$ws.Range("E38").Value = "'  +5.96%  "

$ws.Range("D39").Value = "'11.60"
$ws.Range("E39").Value = "'  +5.43%  "

$ws.Range("B40").Value = "'TheSandbox"
$ws.Range("C40").Value = "'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D40").Value = "'0.6395"
$ws.Range("E40").Value = "'  +6.03%  "

$ws.Range("B41").Value = "'InternetComputer(DFINITY)"
$ws.Range("C41").Value = "'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D41").Value = "'5.028"
$ws.Range("E41").Value = "'  +5.58%  "

$ws.Range("E42").Value = "'  -0.06%  "

$ws.Range("D43").Value = "'1.158"
$ws.Range("E43").Value = "'  +4.25%  "

$ws.Range("D44").Value = "'1.383"
$ws.Range("E44").Value = "'  -2.86%  "

$ws.Range("D45").Value = "'13.52"
$ws.Range("E45").Value = "'  +5.19%  "

$ws.Range("D46").Value = "'0.5980"
$ws.Range("E46").Value = "'  +5.94%  "

$ws.Range("D47").Value = "'3.680"
$ws.Range("E47").Value = "'  +2.48%  "

$ws.Range("D48").Value = "'123.79"
$ws.Range("E48").Value = "'  +3.50%  "

$ws.Range("D49").Value = "'1.974"
$ws.Range("E49").Value = "'  +7.14%  "

$ws.Range("E50").Value = "'  +4.15%  "

$ws.Range("D51").Value = "'0.06889"
$ws.Range("E51").Value = "'  +3.55%  "
